# Insert a new "Password (Mat khau)" list item right after the
# "Luong co ban" bullet, matching its paragraph/run formatting
# (ListParagraph style, ilvl=3, numId=2, Times New Roman 16pt).

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("Lương cơ bản", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target paragraph 'Luong co ban'"
}

$anchorPara = $rng.Paragraphs.Item(1)
$anchorIndex = $anchorPara.Index

# Duplicates the paragraph's pPr/rPr (style, numbering, fonts, size) onto
# a new, empty paragraph inserted immediately after it.
$anchorPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Range.InsertBefore("Password (Mật khẩu)")
